# adding david's additions to spreadsheet
#
# David Kim's user stories are rows 4-5 (Wifi Password Cracker) and 14-15
# (Currency detector). The second line of each story's Must/Should/Could/
# Won't-Have block (rows 13 and 15) was blank; this fills in David's
# additional requirements that were missing from the first commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - second line of requirements for the Wifi Password Cracker story
$ws.Range("E13").Value = "Be able to send password change request or alerts about the network."
$ws.Range("F13").Value = "A GUI to show how strong the wifi connection is."
$ws.Range("G13").Value = "Able to block user from changing password."
$ws.Range("H13").Value = "Mobile phone access"

# Row 15 - second line of requirements for the currency detector app story
$ws.Range("E15").Value = "Able to run the app purely on voice command."
$ws.Range("F15").Value = "Easy to use format"
$ws.Range("G15").Value = "Able to differentiate US currency and other countries currency. "
$ws.Range("H15").Value = "Any functions or features that require payments."

# Leave the selection where the author's last edit landed
$ws.Range("D15").Select() | Out-Null
